$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 4133.16
$ws.Range("J40").Value = 5031.4614
$ws.Range("L40").Value = 5031.4614
$ws.Range("N40").Value = -5381.4614
$ws.Range("H41").Value = 1175.0476
$ws.Range("I41").Value = 906.5
$ws.Range("J41").Value = 1419.1818
$ws.Range("K41").Value = 906.5
$ws.Range("L41").Value = 1419.1818
$ws.Range("M41").Value = -466.5
$ws.Range("N41").Value = -2299.1818
$ws.Range("H76").Value = 8328.286
$ws.Range("I76").Value = 7999
$ws.Range("J76").Value = 8383.166999999999
$ws.Range("K76").Value = 7999
$ws.Range("L76").Value = 8383.166999999999
$ws.Range("M76").Value = -7684
$ws.Range("N76").Value = -9013.166999999999
$ws.Range("H79").Value = 8328.286
$ws.Range("I79").Value = 7999
$ws.Range("J79").Value = 8383.166999999999
$ws.Range("K79").Value = 7999
$ws.Range("L79").Value = 8383.166999999999
$ws.Range("M79").Value = -6907
$ws.Range("N79").Value = -10567.167
$ws.Range("H86").Value = 5266897.5
$ws.Range("I86").Value = 4612.5
$ws.Range("J86").Value = 8775087
$ws.Range("K86").Value = 4612.5
$ws.Range("L86").Value = 8775087
$ws.Range("M86").Value = -3489.5
$ws.Range("N86").Value = -8777333
$ws.Range("H89").Value = 5266897.5
$ws.Range("I89").Value = 4612.5
$ws.Range("J89").Value = 8775087
$ws.Range("K89").Value = 23062.5
$ws.Range("L89").Value = 43875435
$ws.Range("M89").Value = -17446.5
$ws.Range("N89").Value = -43886667
$ws.Range("H94").Value = 1019.7
$ws.Range("I94").Value = 1019.7
$ws.Range("K94").Value = 1019.7
$ws.Range("M94").Value = -568.7
$ws.Range("H98").Value = 1569.7894
$ws.Range("I98").Value = 1642.7059
$ws.Range("J98").Value = 950
$ws.Range("K98").Value = 1642.7059
$ws.Range("L98").Value = 950
$ws.Range("M98").Value = -144.7058999999999
$ws.Range("N98").Value = -3946
$ws.Range("H106").Value = 2950.125
$ws.Range("I106").Value = 2950.125
$ws.Range("K106").Value = 2950.125
$ws.Range("M106").Value = -2319.125
$ws.Range("H122").Value = 1569.7894
$ws.Range("I122").Value = 1642.7059
$ws.Range("J122").Value = 950
$ws.Range("K122").Value = 4928.1177
$ws.Range("L122").Value = 2850
$ws.Range("M122").Value = -2478.1177
$ws.Range("N122").Value = -7750
$ws.Range("H134").Value = 79999.89
$ws.Range("J134").Value = 82500
$ws.Range("L134").Value = 82500
$ws.Range("N134").Value = -92640
$ws.Range("H135").Value = 10004243
$ws.Range("I135").Value = 12195931
$ws.Range("J135").Value = 19887.223
$ws.Range("K135").Value = 109763379
$ws.Range("L135").Value = 178985.007
$ws.Range("M135").Value = -109760844
$ws.Range("N135").Value = -184055.007
$ws.Range("H136").Value = 80000
$ws.Range("I136").Value = 0
$ws.Range("K136").Value = 0
$ws.Range("M136").ClearContents()
$ws.Range("H137").Value = 4844.9473
$ws.Range("I137").Value = 5629.815
$ws.Range("J137").Value = 2918.4546
$ws.Range("K137").Value = 16889.445
$ws.Range("L137").Value = 8755.363799999999
$ws.Range("M137").Value = -14339.445
$ws.Range("N137").Value = -13855.3638
$ws.Range("H138").Value = 4435
$ws.Range("I138").Value = 2898.1428
$ws.Range("J138").Value = 4654.551
$ws.Range("K138").Value = 8694.428400000001
$ws.Range("L138").Value = 13963.653
$ws.Range("M138").Value = -3554.428400000001
$ws.Range("N138").Value = -24243.653
$ws.Range("H139").Value = 100000
$ws.Range("I139").Value = 0
$ws.Range("K139").Value = 0
$ws.Range("M139").ClearContents()
$ws.Range("H140").Value = 0
$ws.Range("J140").Value = 0
$ws.Range("L140").Value = 0
$ws.Range("N140").ClearContents()
$ws.Range("H141").Value = 11027.546
$ws.Range("I141").Value = 6638.8
$ws.Range("J141").Value = 14684.833
$ws.Range("K141").Value = 19916.4
$ws.Range("L141").Value = 44054.499
$ws.Range("M141").Value = -14736.4
$ws.Range("N141").Value = -54414.499

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12504482
$ws.Range("I32").Value = 13893561
$ws.Range("K32").Value = 13893561
$ws.Range("M32").Value = -13893274
$ws.Range("H76").Value = 42652.715
$ws.Range("J76").Value = 45594.832
$ws.Range("L76").Value = 45594.832
$ws.Range("N76").Value = -46270.832
$ws.Range("H79").Value = 42652.715
$ws.Range("J79").Value = 45594.832
$ws.Range("L79").Value = 45594.832
$ws.Range("N79").Value = -47934.832
$ws.Range("H122").Value = 3893.9565
$ws.Range("I122").Value = 2444
$ws.Range("J122").Value = 4405.706
$ws.Range("K122").Value = 7332
$ws.Range("L122").Value = 13217.118
$ws.Range("M122").Value = -4882
$ws.Range("N122").Value = -18117.118

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1649.8334
$ws.Range("I20").Value = 1624.75
$ws.Range("K20").Value = 1624.75
$ws.Range("M20").Value = -1377.75
$ws.Range("H105").Value = 4004.25
$ws.Range("I105").Value = 2869.9092
$ws.Range("K105").Value = 2869.9092
$ws.Range("M105").Value = -1122.9092
$ws.Range("H127").Value = 48000
$ws.Range("J127").Value = 48000
$ws.Range("L127").Value = 48000
$ws.Range("N127").Value = -57920

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1831.7778
$ws.Range("I31").Value = 1542.7858
$ws.Range("K31").Value = 1542.7858
$ws.Range("M31").Value = -1247.7858
$ws.Range("H34").Value = 1831.7778
$ws.Range("I34").Value = 1542.7858
$ws.Range("K34").Value = 1542.7858
$ws.Range("M34").Value = -1340.7858
$ws.Range("H99").Value = 16032742
$ws.Range("J99").Value = 20006024
$ws.Range("L99").Value = 20006024
$ws.Range("N99").Value = -20009020
$ws.Range("H107").Value = 11442.789
$ws.Range("I107").Value = 914.9
$ws.Range("J107").Value = 23140.445
$ws.Range("K107").Value = 914.9
$ws.Range("L107").Value = 23140.445
$ws.Range("M107").Value = 1005.1
$ws.Range("N107").Value = -26980.445
$ws.Range("H122").Value = 572219.75
$ws.Range("I122").Value = 1276996
$ws.Range("K122").Value = 3830988
$ws.Range("M122").Value = -3828538
$ws.Range("H126").Value = 16032742
$ws.Range("J126").Value = 20006024
$ws.Range("L126").Value = 60018072
$ws.Range("N126").Value = -60023012

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H92").Value = 611.1111
$ws.Range("J92").Value = 875
$ws.Range("L92").Value = 2625
$ws.Range("N92").Value = -5121
$ws.Range("H122").Value = 805.4706
$ws.Range("I122").Value = 858.5454999999999
$ws.Range("K122").Value = 7726.9095
$ws.Range("M122").Value = -5276.9095

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H55").Value = 4689.8887
$ws.Range("J55").Value = 5284
$ws.Range("L55").Value = 5284
$ws.Range("N55").Value = -5938
$ws.Range("H70").Value = 149525.12
$ws.Range("I70").Value = 284053
$ws.Range("J70").Value = 14997.25
$ws.Range("K70").Value = 284053
$ws.Range("L70").Value = 14997.25
$ws.Range("M70").Value = -283783
$ws.Range("N70").Value = -15537.25
$ws.Range("H73").Value = 149525.12
$ws.Range("I73").Value = 284053
$ws.Range("J73").Value = 14997.25
$ws.Range("K73").Value = 284053
$ws.Range("L73").Value = 14997.25
$ws.Range("M73").Value = -283117
$ws.Range("N73").Value = -16869.25
$ws.Range("H80").Value = 33447336
$ws.Range("J80").Value = 41669656
$ws.Range("L80").Value = 41669656
$ws.Range("N80").Value = -41671652
$ws.Range("H83").Value = 33447336
$ws.Range("J83").Value = 41669656
$ws.Range("L83").Value = 208348280
$ws.Range("N83").Value = -208358264
$ws.Range("H122").Value = 5775.1055
$ws.Range("I122").Value = 4613.3335
$ws.Range("J122").Value = 6820.7
$ws.Range("K122").Value = 13840.0005
$ws.Range("L122").Value = 20462.1
$ws.Range("M122").Value = -11390.0005
$ws.Range("N122").Value = -25362.1

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H35").Value = 10267.5
$ws.Range("I35").Value = 500
$ws.Range("J35").Value = 20035
$ws.Range("K35").Value = 500
$ws.Range("L35").Value = 20035
$ws.Range("M35").Value = -164
$ws.Range("N35").Value = -20707
$ws.Range("H68").Value = 3997.5454
$ws.Range("I68").Value = 3994.6
$ws.Range("K68").Value = 3994.6
$ws.Range("M68").Value = -3245.6
$ws.Range("H71").Value = 3997.5454
$ws.Range("I71").Value = 3994.6
$ws.Range("K71").Value = 19973
$ws.Range("M71").Value = -16229
$ws.Range("H100").Value = 283002.75
$ws.Range("I100").Value = 558005.5
$ws.Range("K100").Value = 558005.5
$ws.Range("M100").Value = -557464.5
$ws.Range("H109").Value = 37352.53
$ws.Range("J109").Value = 37352.53
$ws.Range("L109").Value = 37352.53
$ws.Range("N109").Value = -40126.53
$ws.Range("H136").Value = 3463.6216
$ws.Range("I136").Value = 3447.0312
$ws.Range("J136").Value = 3569.8
$ws.Range("K136").Value = 10341.0936
$ws.Range("L136").Value = 10709.4
$ws.Range("M136").Value = -7791.0936
$ws.Range("N136").Value = -15809.4

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H32").Value = 15756.5
$ws.Range("I32").Value = 4342
$ws.Range("K32").Value = 4342
$ws.Range("M32").Value = -4025
$ws.Range("H96").Value = 108200.6
$ws.Range("I96").Value = 254000
$ws.Range("J96").Value = 11001
$ws.Range("K96").Value = 254000
$ws.Range("L96").Value = 11001
$ws.Range("M96").Value = -252627
$ws.Range("N96").Value = -13747
$ws.Range("H122").Value = 2822.1
$ws.Range("I122").Value = 2865.9524
$ws.Range("J122").Value = 2719.7778
$ws.Range("K122").Value = 8597.8572
$ws.Range("L122").Value = 8159.3334
$ws.Range("M122").Value = -6147.8572
$ws.Range("N122").Value = -13059.3334
$ws.Range("H132").Value = 2320.484
$ws.Range("I132").Value = 1622.5769
$ws.Range("J132").Value = 5949.6
$ws.Range("K132").Value = 4867.7307
$ws.Range("L132").Value = 17848.8
$ws.Range("M132").Value = -2337.7307
$ws.Range("N132").Value = -22908.8
